$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1554.6666
$ws.Range("I19").Value = 1698.6
$ws.Range("K19").Value = 1698.6
$ws.Range("M19").Value = -1523.6
$ws.Range("H68").Value = 10000
$ws.Range("I68").Value = 10000
$ws.Range("K68").Value = 10000
$ws.Range("M68").Value = -9251
$ws.Range("H71").Value = 10000
$ws.Range("I71").Value = 10000
$ws.Range("K71").Value = 30000
$ws.Range("M71").Value = -26256
$ws.Range("H98").Value = 1540.2858
$ws.Range("J98").Value = 1899.5
$ws.Range("L98").Value = 1899.5
$ws.Range("N98").Value = -4895.5
$ws.Range("H106").Value = 11983.417
$ws.Range("I106").Value = 11983.417
$ws.Range("K106").Value = 11983.417
$ws.Range("M106").Value = -11352.417
$ws.Range("H112").Value = 2653.4285
$ws.Range("J112").Value = 3119.818
$ws.Range("L112").Value = 9359.454000000002
$ws.Range("N112").Value = -11575.454
$ws.Range("H122").Value = 1540.2858
$ws.Range("J122").Value = 1899.5
$ws.Range("L122").Value = 5698.5
$ws.Range("N122").Value = -10598.5
$ws.Range("H137").Value = 10102445
$ws.Range("I137").Value = 15873961
$ws.Range("J137").Value = 2291.8333
$ws.Range("K137").Value = 47621883
$ws.Range("L137").Value = 6875.499899999999
$ws.Range("M137").Value = -47619333
$ws.Range("N137").Value = -11975.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1460.4286
$ws.Range("I45").Value = 1287.1666
$ws.Range("K45").Value = 1287.1666
$ws.Range("M45").Value = -910.1666
$ws.Range("H92").Value = 189500
$ws.Range("J92").Value = 189500
$ws.Range("L92").Value = 189500
$ws.Range("N92").Value = -194492
$ws.Range("H132").Value = 4959.8887
$ws.Range("I132").Value = 5097.8823
$ws.Range("J132").Value = 2614
$ws.Range("K132").Value = 15293.6469
$ws.Range("L132").Value = 7842
$ws.Range("M132").Value = -12763.6469
$ws.Range("N132").Value = -12902

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 602.3333
$ws.Range("I22").Value = 641.5
$ws.Range("K22").Value = 641.5
$ws.Range("M22").Value = -468.5
$ws.Range("H86").Value = 28538.035
$ws.Range("I86").Value = 28557.777
$ws.Range("J86").Value = 28502.5
$ws.Range("K86").Value = 28557.777
$ws.Range("L86").Value = 28502.5
$ws.Range("M86").Value = -27434.777
$ws.Range("N86").Value = -30748.5
$ws.Range("H89").Value = 28538.035
$ws.Range("I89").Value = 28557.777
$ws.Range("J89").Value = 28502.5
$ws.Range("K89").Value = 142788.885
$ws.Range("L89").Value = 142512.5
$ws.Range("M89").Value = -137172.885
$ws.Range("N89").Value = -153744.5
$ws.Range("H99").Value = 1011
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2666.8333
$ws.Range("I6").Value = 1859.8
$ws.Range("K6").Value = 1859.8
$ws.Range("M6").Value = -1746.8
$ws.Range("H7").Value = 62500204
$ws.Range("J7").Value = 14
$ws.Range("L7").Value = 14
$ws.Range("N7").Value = -240
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -20952
$ws.Range("H58").Value = 2454.2104
$ws.Range("I58").Value = 2620.6875
$ws.Range("K58").Value = 2620.6875
$ws.Range("M58").Value = -2417.6875
$ws.Range("H59").Value = 148000
$ws.Range("J59").Value = 148000
$ws.Range("L59").Value = 148000
$ws.Range("N59").Value = -150290
$ws.Range("H99").Value = 4494.9
$ws.Range("I99").Value = 4010.6
$ws.Range("K99").Value = 4010.6
$ws.Range("M99").Value = -2512.6
$ws.Range("H126").Value = 4494.9
$ws.Range("I126").Value = 4010.6
$ws.Range("K126").Value = 12031.8
$ws.Range("M126").Value = -9561.799999999999
$ws.Range("H134").Value = 2197.76
$ws.Range("J134").Value = 1915
$ws.Range("L134").Value = 5745
$ws.Range("N134").Value = -10815
$ws.Range("H136").Value = 2454.2104
$ws.Range("I136").Value = 2620.6875
$ws.Range("K136").Value = 7862.0625
$ws.Range("M136").Value = -5312.0625
$ws.Range("H141").Value = 349351.22
$ws.Range("J141").Value = 386770.12
$ws.Range("L141").Value = 386770.12
$ws.Range("N141").Value = -397130.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 42
$ws.Range("I38").Value = 38.77778
$ws.Range("K38").Value = 116.33334
$ws.Range("M38").Value = 230.66666
$ws.Range("H113").Value = 2490.6365
$ws.Range("J113").Value = 2518.4
$ws.Range("L113").Value = 7555.200000000001
$ws.Range("N113").Value = -11895.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("H102").Value = 1854.7894
$ws.Range("I102").Value = 1899.4706
$ws.Range("J102").Value = 1475
$ws.Range("K102").Value = 1899.4706
$ws.Range("L102").Value = 1475
$ws.Range("M102").Value = -277.4706000000001
$ws.Range("N102").Value = -4719
$ws.Range("M53").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2947.8572
$ws.Range("I122").Value = 2736.5
$ws.Range("J122").Value = 3229.6667
$ws.Range("K122").Value = 8209.5
$ws.Range("L122").Value = 9689.000100000001
$ws.Range("M122").Value = -5759.5
$ws.Range("N122").Value = -14589.0001
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9470
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10586
$ws.Range("H81").Value = 8513
$ws.Range("I81").Value = 4298.9
$ws.Range("J81").Value = 13780.625
$ws.Range("K81").Value = 8597.799999999999
$ws.Range("L81").Value = 27561.25
$ws.Range("M81").Value = -7536.799999999999
$ws.Range("N81").Value = -29683.25
$ws.Range("H82").Value = 70273
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H84").Value = 8513
$ws.Range("I84").Value = 4298.9
$ws.Range("J84").Value = 13780.625
$ws.Range("K84").Value = 42989
$ws.Range("L84").Value = 137806.25
$ws.Range("M84").Value = -37685
$ws.Range("N84").Value = -148414.25
$ws.Range("H85").Value = 70273
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H122").Value = 1521.1111
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H132").Value = 58837020
$ws.Range("I132").Value = 15089.333
$ws.Range("K132").Value = 45267.999
$ws.Range("M132").Value = -42737.999
$ws.Range("M26").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("N122").ClearContents()
